# Apply the Betfair Back/Lay odds update for 2025-12-09.
# - Rows 2-13: individual odds values changed (per the diff).
# - Row 14: a brand-new match (Bolivian Liga de Futbol Profesional) appended.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing odds values on rows 2-13 ---
$ws.Range("G2").Value = 1.7
$ws.Range("I2").Value = 9.6
$ws.Range("N2").Value = 3.1
$ws.Range("Q2").Value = 1.92
$ws.Range("V2").Value = 1.11
$ws.Range("W2").Value = 2.42
$ws.Range("F3").Value = 1.81
$ws.Range("I3").Value = 6.6
$ws.Range("J3").Value = 2.96
$ws.Range("M3").Value = 1.09
$ws.Range("F4").Value = 1.78
$ws.Range("I4").Value = 6.4
$ws.Range("J4").Value = 3.25
$ws.Range("K4").Value = 3.65
$ws.Range("L4").Value = 1.54
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 2.64
$ws.Range("O4").Value = 1.5
$ws.Range("P4").Value = 1.55
$ws.Range("Q4").Value = 2.46
$ws.Range("R4").Value = 1.2
$ws.Range("S4").Value = 5
$ws.Range("T4").Value = 2.2
$ws.Range("U4").Value = 1.68
$ws.Range("V4").Value = 1.19
$ws.Range("X4").Value = 11
$ws.Range("AA4").Value = 210
$ws.Range("AB4").Value = 7.4
$ws.Range("AF4").Value = 11.5
$ws.Range("AI4").Value = 160
$ws.Range("G5").Value = 7.8
$ws.Range("H5").Value = 1.53
$ws.Range("I5").Value = 1.54
$ws.Range("J5").Value = 4.5
$ws.Range("L5").Value = 1.36
$ws.Range("R5").Value = 1.47
$ws.Range("S5").Value = 3
$ws.Range("U5").Value = 1.98
$ws.Range("V5").Value = 2.86
$ws.Range("X5").Value = 16.5
$ws.Range("Z5").Value = 8.800000000000001
$ws.Range("AL5").Value = 95
$ws.Range("AN5").Value = 130
$ws.Range("AO5").Value = 7.2
$ws.Range("F6").Value = 1.24
$ws.Range("H6").Value = 13.5
$ws.Range("I6").Value = 14.5
$ws.Range("J6").Value = 7.6
$ws.Range("L6").Value = 1.2
$ws.Range("N6").Value = 8.4
$ws.Range("O6").Value = 1.12
$ws.Range("P6").Value = 3.45
$ws.Range("Q6").Value = 1.38
$ws.Range("R6").Value = 2
$ws.Range("T6").Value = 1.9
$ws.Range("U6").Value = 2.04
$ws.Range("W6").Value = 5
$ws.Range("Y6").Value = 710
$ws.Range("AA6").Value = 580
$ws.Range("AB6").Value = 15
$ws.Range("AC6").Value = 17.5
$ws.Range("AD6").Value = 48
$ws.Range("AE6").Value = 190
$ws.Range("AF6").Value = 9.800000000000001
$ws.Range("AG6").Value = 12
$ws.Range("AH6").Value = 30
$ws.Range("AJ6").Value = 9.6
$ws.Range("AL6").Value = 34
$ws.Range("AM6").Value = 130
$ws.Range("AN6").Value = 3.2
$ws.Range("AO6").Value = 150
$ws.Range("H7").Value = 2.28
$ws.Range("I7").Value = 2.32
$ws.Range("J7").Value = 3.65
$ws.Range("K7").Value = 3.75
$ws.Range("N7").Value = 4.7
$ws.Range("O7").Value = 1.25
$ws.Range("R7").Value = 1.5
$ws.Range("T7").Value = 1.63
$ws.Range("U7").Value = 2.5
$ws.Range("Z7").Value = 15.5
$ws.Range("AE7").Value = 21
$ws.Range("AL7").Value = 40
$ws.Range("AM7").Value = 65
$ws.Range("F8").Value = 1.44
$ws.Range("H8").Value = 8.4
$ws.Range("I8").Value = 8.6
$ws.Range("J8").Value = 5.2
$ws.Range("K8").Value = 5.3
$ws.Range("L8").Value = 1.36
$ws.Range("W8").Value = 3.2
$ws.Range("X8").Value = 19
$ws.Range("AA8").Value = 310
$ws.Range("AE8").Value = 140
$ws.Range("F9").Value = 3.15
$ws.Range("G9").Value = 3.2
$ws.Range("H9").Value = 2.3
$ws.Range("J9").Value = 3.95
$ws.Range("O9").Value = 1.18
$ws.Range("P9").Value = 2.66
$ws.Range("R9").Value = 1.68
$ws.Range("S9").Value = 2.4
$ws.Range("U9").Value = 2.8
$ws.Range("W9").Value = 1.45
$ws.Range("AF9").Value = 26
$ws.Range("AL9").Value = 32
$ws.Range("AN9").Value = 18
$ws.Range("AO9").Value = 11
$ws.Range("F10").Value = 2.32
$ws.Range("G10").Value = 2.34
$ws.Range("J10").Value = 3.95
$ws.Range("P10").Value = 2.64
$ws.Range("AA10").Value = 55
$ws.Range("AC10").Value = 9.4
$ws.Range("AO10").Value = 18
$ws.Range("F11").Value = 2.2
$ws.Range("G11").Value = 2.22
$ws.Range("H11").Value = 3.5
$ws.Range("I11").Value = 3.55
$ws.Range("K11").Value = 3.8
$ws.Range("P11").Value = 2.3
$ws.Range("S11").Value = 2.8
$ws.Range("V11").Value = 1.38
$ws.Range("W11").Value = 1.82
$ws.Range("AD11").Value = 14
$ws.Range("H12").Value = 19.5
$ws.Range("I12").Value = 20
$ws.Range("L12").Value = 1.14
$ws.Range("N12").Value = 11.5
$ws.Range("P12").Value = 4.6
$ws.Range("R12").Value = 2.44
$ws.Range("S12").Value = 1.66
$ws.Range("W12").Value = 6.6
$ws.Range("X12").Value = 940
$ws.Range("Z12").Value = 260
$ws.Range("AB12").Value = 19
$ws.Range("AC12").Value = 25
$ws.Range("AE12").Value = 280
$ws.Range("AG12").Value = 14
$ws.Range("AH12").Value = 36
$ws.Range("AJ12").Value = 11
$ws.Range("AN12").Value = 2.42
$ws.Range("G13").Value = 3.15
$ws.Range("H13").Value = 2.4
$ws.Range("I13").Value = 2.44
$ws.Range("P13").Value = 2.3
$ws.Range("Q13").Value = 1.72
$ws.Range("S13").Value = 2.8
$ws.Range("T13").Value = 1.62
$ws.Range("U13").Value = 2.52
$ws.Range("V13").Value = 1.69
$ws.Range("W13").Value = 1.47
$ws.Range("AO13").Value = 15

# --- Append new row 14: Jorge Wilstermann vs CD Gualberto Villarroel ---
$ws.Range("A14").Value = "Bolivian Liga de Futbol Profesional"

# Date column (B) must stay literal text "2025-12-09", like the rest of column B,
# rather than Excel auto-converting it into a date serial number.
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "2025-12-09"
$ws.Range("B14").Style = "Normal"

$ws.Range("C14").Value = "20:00:00"
$ws.Range("D14").Value = "Jorge Wilstermann"
$ws.Range("E14").Value = "CD Gualberto Villarroel"

$ws.Range("F14").Value = 2.42
$ws.Range("G14").Value = 2.72
$ws.Range("H14").Value = 2.72
$ws.Range("I14").Value = 3.05
$ws.Range("J14").Value = 3.35
$ws.Range("K14").Value = 4
$ws.Range("L14").Value = 1.35
$ws.Range("M14").Value = 1.01
$ws.Range("N14").Value = 3.95
$ws.Range("O14").Value = 1.26
$ws.Range("P14").Value = 2.08
$ws.Range("Q14").Value = 1.77
$ws.Range("R14").Value = 1.41
$ws.Range("S14").Value = 2.94
$ws.Range("T14").Value = 1.64
$ws.Range("U14").Value = 2.24
$ws.Range("V14").Value = 1.5
$ws.Range("W14").Value = 1.58
$ws.Range("X14").Value = 23
$ws.Range("Y14").Value = 16.5
$ws.Range("Z14").Value = 26
$ws.Range("AA14").Value = 55
$ws.Range("AB14").Value = 15.5
$ws.Range("AC14").Value = 11
$ws.Range("AD14").Value = 15.5
$ws.Range("AE14").Value = 36
$ws.Range("AF14").Value = 23
$ws.Range("AG14").Value = 15
$ws.Range("AH14").Value = 19.5
$ws.Range("AI14").Value = 46
$ws.Range("AJ14").Value = 44
$ws.Range("AK14").Value = 32
$ws.Range("AL14").Value = 42
$ws.Range("AM14").Value = 90
$ws.Range("AN14").Value = 22
$ws.Range("AO14").Value = 27
